$d = $word.ActiveDocument
$d.Content.Find.Execute("697×2=1394", $true, $false, $false, $false, $false, $true, 1, $false, "565×6=3390", 2) | Out-Null
$d.Content.Find.Execute("155×4=620", $true, $false, $false, $false, $false, $true, 1, $false, "565×2=1130", 2) | Out-Null
$d.Content.Find.Execute("424×4=1696", $true, $false, $false, $false, $false, $true, 1, $false, "670×8=5360", 2) | Out-Null
$d.Content.Find.Execute("213×5=1065", $true, $false, $false, $false, $false, $true, 1, $false, "839×5=4195", 2) | Out-Null
$d.Content.Find.Execute("836×4=3344", $true, $false, $false, $false, $false, $true, 1, $false, "563×6=3378", 2) | Out-Null
$d.Content.Find.Execute("700×5=3500", $true, $false, $false, $false, $false, $true, 1, $false, "849×4=3396", 2) | Out-Null
$d.Content.Find.Execute("294×5=1470", $true, $false, $false, $false, $false, $true, 1, $false, "543×4=2172", 2) | Out-Null
$d.Content.Find.Execute("498×2=996", $true, $false, $false, $false, $false, $true, 1, $false, "194×2=388", 2) | Out-Null
$d.Content.Find.Execute("270×7=1890", $true, $false, $false, $false, $false, $true, 1, $false, "513×5=2565", 2) | Out-Null
$d.Content.Find.Execute("355×5=1775", $true, $false, $false, $false, $false, $true, 1, $false, "730×4=2920", 2) | Out-Null
$d.Content.Find.Execute("731×3=2193", $true, $false, $false, $false, $false, $true, 1, $false, "728×4=2912", 2) | Out-Null
$d.Content.Find.Execute("340×9=3060", $true, $false, $false, $false, $false, $true, 1, $false, "537×3=1611", 2) | Out-Null
$d.Content.Find.Execute("920×6=5520", $true, $false, $false, $false, $false, $true, 1, $false, "803×8=6424", 2) | Out-Null
$d.Content.Find.Execute("113×8=904", $true, $false, $false, $false, $false, $true, 1, $false, "694×9=6246", 2) | Out-Null
$d.Content.Find.Execute("803×7=5621", $true, $false, $false, $false, $false, $true, 1, $false, "464×2=928", 2) | Out-Null
$d.Content.Find.Execute("314×7=2198", $true, $false, $false, $false, $false, $true, 1, $false, "572×5=2860", 2) | Out-Null
$d.Content.Find.Execute("638×6=3828", $true, $false, $false, $false, $false, $true, 1, $false, "759×3=2277", 2) | Out-Null
$d.Content.Find.Execute("778×9=7002", $true, $false, $false, $false, $false, $true, 1, $false, "653×3=1959", 2) | Out-Null
$d.Content.Find.Execute("386×4=1544", $true, $false, $false, $false, $false, $true, 1, $false, "554×9=4986", 2) | Out-Null
$d.Content.Find.Execute("318×5=1590", $true, $false, $false, $false, $false, $true, 1, $false, "534×6=3204", 2) | Out-Null
$d.Content.Find.Execute("726×9=6534", $true, $false, $false, $false, $false, $true, 1, $false, "997×7=6979", 2) | Out-Null
$d.Content.Find.Execute("468×4=1872", $true, $false, $false, $false, $false, $true, 1, $false, "923×7=6461", 2) | Out-Null
$d.Content.Find.Execute("887×6=5322", $true, $false, $false, $false, $false, $true, 1, $false, "562×3=1686", 2) | Out-Null
$d.Content.Find.Execute("500×4=2000", $true, $false, $false, $false, $false, $true, 1, $false, "248×8=1984", 2) | Out-Null
$d.Content.Find.Execute("481×5=2405", $true, $false, $false, $false, $false, $true, 1, $false, "396×4=1584", 2) | Out-Null
